$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Row 3 on every sheet corresponds to the
# 5a765ba2-0408-43b4-afa6-c5f520b8e3dd file, which is now ready for handoff.

# Overview sheet: Status columns for zh-cn (B) and de-de (C), plus the
# Latest Handoff Date (D) for that file.
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"
$overview.Range("D3").Value = "2016-03-21 03:23:47"

# zh-cn sheet: Status (C) and Latest Handoff Datetime (E)
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("E3").Value = "2016-03-21 03:23:39"

# de-de sheet: Status (C) and Latest Handoff Datetime (E)
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("E3").Value = "2016-03-21 03:23:47"
